# Applies the Tue Jan 3 04:30:00 UTC 2023 GitHub Actions "symbol list" refresh:
# coin rows shift by one position (new row inserted at the top of the shuffled
# block) and Price / Volume(1h) figures are refreshed across the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / Link URL columns - plain text, never numeric-looking, so a direct
# assignment is safe and will not be coerced to a number by Excel.
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'

# Price / Volume(1h) columns look numeric ("245.58", "1.07%", "1,597.14%", ...)
# but must stay literal text (matching the original inline-string cells) instead
# of being parsed into Excel numbers/percentages. Force the cell to Text format
# before assigning, then restore the default "Normal" style so the only change
# that sticks is the cell value itself.
$numericLooking = [ordered]@{
    "D2" = '245.58'
    "E2" = '1.07%'
    "D3" = '29.39'
    "E3" = '-2.41%'
    "D4" = '5.148'
    "E4" = '0.07%'
    "D5" = '0.05776'
    "E5" = '1.84%'
    "D6" = '6.612'
    "E6" = '1.41%'
    "D7" = '3.165'
    "E7" = '4.87%'
    "D8" = '0.8573'
    "E8" = '2.12%'
    "D9" = '0.8565'
    "E9" = '-0.37%'
    "D10" = '0.01021'
    "E10" = '1,597.14%'
    "D11" = '0.1363'
    "E11" = '1.99%'
    "D12" = '0.07018'
    "E12" = '1.51%'
    "D13" = '0.02985'
    "E13" = '4.32%'
    "D14" = '0.09359'
    "E14" = '-0.12%'
    "D15" = '0.001547'
    "E15" = '1.74%'
    "D16" = '0.005986'
    "E16" = '-0.20%'
    "D17" = '3.482'
    "E17" = '-1.16%'
    "D18" = '2.170'
    "E18" = '2.05%'
    "D19" = '0.3203'
    "E19" = '1.65%'
    "D20" = '0.03314'
    "E20" = '1.72%'
    "D21" = '0.1282'
    "E21" = '-1.02%'
    "D22" = '3.317'
    "E22" = '-7.94%'
    "D23" = '0.04134'
    "E23" = '-0.52%'
    "E24" = '1.94%'
    "E25" = '1.23%'
    "E26" = '-7.07%'
    "E27" = '2.61%'
    "D28" = '0.0001445'
    "E28" = '3.44%'
    "D40" = '0.03728'
    "E40" = '0.42%'
    "D41" = '0.005900'
    "E41" = '12.88%'
    "D42" = '0.1070'
    "E42" = '1.29%'
    "D43" = '0.002200'
    "E43" = '-4.75%'
    "D44" = '0.008539'
    "E44" = '-12.53%'
    "D45" = '0.00005289'
    "E45" = '3.60%'
    "E46" = '0.06%'
    "D47" = '0.05800'
    "E47" = '-41.97%'
    "E48" = '-19.60%'
    "D49" = '0.00002100'
    "E49" = '0.06%'
    "D50" = '0.0002000'
    "E50" = '0.06%'
}
foreach ($addr in $numericLooking.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLooking[$addr]
    $cell.Style = "Normal"
}
